$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 170, pushing existing rows 170-176 down to 171-177.
$ws.Rows(170).Insert()

# Populate the newly inserted row 170 with the new weekly price-report entry.
$ws.Range("A170").Value = 10
$ws.Range("B170").Value = "Vega Modelo de Temuco"
$ws.Range("C170").Value = "La Araucanía"
$ws.Range("D170").Value = 44516
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 100112039
$ws.Range("G170").Value = "Ciboulette"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 65
$ws.Range("K170").Value = 5000
$ws.Range("L170").Value = 5000
$ws.Range("M170").Value = 5000
$ws.Range("N170").Value = "$/docena de atados"
$ws.Range("O170").Value = "Provincia de Cautín"
$ws.Range("P170").Value = 1667
$ws.Range("Q170").Value = 3
$ws.Range("R170").Value = "Hortaliza"
